# Auto-generated COM script reproducing the PlayerPerformance_3749.xlsx restructuring:
#   - insert a new "Player Info" sheet at the front
#   - keep "ODI Batting" / "ODI Bowling" sheets, renaming MATCH_CARD_LINK -> MATCH_CODE
#     and collapsing the scorecard URL down to the bare numeric match code
#   - append a new "ODI Batting Extra" sheet at the end
$wb = $excel.ActiveWorkbook

# --- "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE, URL -> bare match code, clear blank INNING_NUMBER cells ---
$wsBatting = $wb.Worksheets.Item("ODI Batting")
$wsBatting.Cells.Item(1, 4).Value = "MATCH_CODE"
$wsBatting.Cells.Item(2, 4).Value = "3068"
$wsBatting.Cells.Item(3, 4).Value = "3070"
$wsBatting.Cells.Item(4, 4).Value = "3119"
$wsBatting.Cells.Item(5, 4).Value = "3364"
$wsBatting.Cells.Item(6, 4).Value = "3365"
$wsBatting.Cells.Item(7, 4).Value = "3369"
$wsBatting.Cells.Item(8, 4).Value = "3371"
$wsBatting.Cells.Item(9, 4).Value = "3377"
$wsBatting.Cells.Item(10, 4).Value = "3382"
$wsBatting.Cells.Item(11, 4).Value = "3384"
$wsBatting.Cells.Item(12, 4).Value = "3387"
$wsBatting.Cells.Item(13, 4).Value = "3389"
$wsBatting.Cells.Item(14, 4).Value = "3391"
$wsBatting.Cells.Item(15, 4).Value = "3392"
$wsBatting.Cells.Item(16, 4).Value = "3393"
$wsBatting.Cells.Item(17, 4).Value = "3395"
$wsBatting.Cells.Item(18, 4).Value = "3396"
$wsBatting.Cells.Item(19, 4).Value = "3401"
$wsBatting.Cells.Item(20, 4).Value = "3408"
$wsBatting.Cells.Item(21, 4).Value = "3409"
$wsBatting.Cells.Item(22, 4).Value = "3410"
$wsBatting.Cells.Item(23, 4).Value = "3411"
$wsBatting.Cells.Item(24, 4).Value = "3413"
$wsBatting.Cells.Item(25, 4).Value = "3431"
$wsBatting.Cells.Item(26, 4).Value = "3432"
$wsBatting.Cells.Item(27, 4).Value = "3434"
$wsBatting.Cells.Item(28, 4).Value = "3435"
$wsBatting.Cells.Item(29, 4).Value = "3446"
$wsBatting.Cells.Item(30, 4).Value = "3447"
$wsBatting.Cells.Item(31, 4).Value = "3448"
$wsBatting.Cells.Item(32, 4).Value = "3449"
$wsBatting.Cells.Item(33, 4).Value = "3458"
$wsBatting.Cells.Item(34, 4).Value = "3460"
$wsBatting.Cells.Item(35, 4).Value = "3462"
$wsBatting.Cells.Item(36, 4).Value = "3465"
$wsBatting.Cells.Item(37, 4).Value = "3467"
$wsBatting.Cells.Item(38, 4).Value = "3490"
$wsBatting.Cells.Item(39, 4).Value = "3492"
$wsBatting.Cells.Item(40, 4).Value = "3493"
$wsBatting.Cells.Item(41, 4).Value = "3508"
$wsBatting.Cells.Item(42, 4).Value = "3512"
$wsBatting.Cells.Item(43, 4).Value = "3516"
$wsBatting.Cells.Item(44, 4).Value = "3518"
$wsBatting.Cells.Item(45, 4).Value = "3520"
$wsBatting.Cells.Item(46, 4).Value = "3524"
$wsBatting.Cells.Item(47, 4).Value = "3527"
$wsBatting.Cells.Item(48, 4).Value = "3529"
$wsBatting.Cells.Item(49, 4).Value = "3530"
$wsBatting.Cells.Item(50, 4).Value = "3534"
$wsBatting.Cells.Item(51, 4).Value = "3536"
$wsBatting.Cells.Item(52, 4).Value = "3542"
$wsBatting.Cells.Item(53, 4).Value = "3543"
$wsBatting.Cells.Item(54, 4).Value = "3576"
$wsBatting.Cells.Item(55, 4).Value = "3578"
$wsBatting.Cells.Item(56, 4).Value = "3579"
$wsBatting.Cells.Item(57, 4).Value = "3589"
$wsBatting.Cells.Item(58, 4).Value = "3616"
$wsBatting.Cells.Item(59, 4).Value = "3618"
$wsBatting.Cells.Item(60, 4).Value = "3621"
$wsBatting.Cells.Item(61, 4).Value = "3626"
$wsBatting.Cells.Item(62, 4).Value = "3630"
$wsBatting.Cells.Item(63, 4).Value = "3631"
$wsBatting.Cells.Item(64, 4).Value = "3635"
$wsBatting.Cells.Item(65, 4).Value = "3638"
$wsBatting.Cells.Item(66, 4).Value = "3639"
$wsBatting.Cells.Item(67, 4).Value = "3640"
$wsBatting.Cells.Item(68, 4).Value = "3641"
$wsBatting.Cells.Item(69, 4).Value = "3642"
$wsBatting.Cells.Item(70, 4).Value = "3646"
$wsBatting.Cells.Item(71, 4).Value = "3647"
$wsBatting.Cells.Item(72, 4).Value = "3648"
$wsBatting.Cells.Item(73, 4).Value = "3693"
$wsBatting.Cells.Item(74, 4).Value = "3696"
$wsBatting.Cells.Item(75, 4).Value = "3703"
$wsBatting.Cells.Item(76, 4).Value = "3706"
$wsBatting.Cells.Item(77, 4).Value = "3710"
$wsBatting.Cells.Item(78, 4).Value = "3712"
$wsBatting.Cells.Item(79, 4).Value = "3714"
$wsBatting.Cells.Item(80, 4).Value = "3716"
$wsBatting.Cells.Item(81, 4).Value = "3718"
$wsBatting.Cells.Item(82, 4).Value = "3723"
$wsBatting.Cells.Item(83, 4).Value = "3726"
$wsBatting.Cells.Item(84, 4).Value = "3729"
$wsBatting.Cells.Item(85, 4).Value = "3734"
$wsBatting.Cells.Item(86, 4).Value = "3737"
$wsBatting.Cells.Item(87, 4).Value = "3739"
$wsBatting.Cells.Item(88, 4).Value = "3743"
$wsBatting.Cells.Item(89, 4).Value = "3748"
$wsBatting.Cells.Item(90, 4).Value = "3759"
$wsBatting.Cells.Item(91, 4).Value = "3765"
$wsBatting.Cells.Item(92, 4).Value = "3769"
$wsBatting.Cells.Item(93, 4).Value = "3779"
$wsBatting.Cells.Item(94, 4).Value = "3782"
$wsBatting.Cells.Item(95, 4).Value = "3790"
$wsBatting.Cells.Item(96, 4).Value = "3814"
$wsBatting.Cells.Item(97, 4).Value = "3819"
$wsBatting.Cells.Item(98, 4).Value = "3820"
$wsBatting.Cells.Item(99, 4).Value = "3821"
$wsBatting.Cells.Item(100, 4).Value = "3822"
$wsBatting.Cells.Item(101, 4).Value = "3852"
$wsBatting.Cells.Item(102, 4).Value = "3853"
$wsBatting.Cells.Item(103, 4).Value = "3855"
$wsBatting.Cells.Item(104, 4).Value = "3865"
$wsBatting.Cells.Item(105, 4).Value = "3866"
$wsBatting.Cells.Item(106, 4).Value = "3868"
$wsBatting.Cells.Item(107, 4).Value = "3870"
$wsBatting.Cells.Item(108, 4).Value = "3872"
$wsBatting.Cells.Item(109, 4).Value = "4063"
$wsBatting.Cells.Item(110, 4).Value = "4064"
$wsBatting.Cells.Item(111, 4).Value = "4065"
$wsBatting.Cells.Item(112, 4).Value = "4079"
$wsBatting.Cells.Item(113, 4).Value = "4081"
$wsBatting.Cells.Item(114, 4).Value = "4082"
$wsBatting.Cells.Item(115, 4).Value = "4084"
$wsBatting.Cells.Item(116, 4).Value = "4087"
$wsBatting.Cells.Item(117, 4).Value = "4096"
$wsBatting.Cells.Item(118, 4).Value = "4098"
$wsBatting.Cells.Item(119, 4).Value = "4302"
$wsBatting.Cells.Item(120, 4).Value = "4305"
$wsBatting.Cells.Item(121, 4).Value = "4309"
$wsBatting.Cells.Item(122, 4).Value = "4322"
$wsBatting.Cells.Item(123, 4).Value = "4344"
$wsBatting.Cells.Item(124, 4).Value = "4350"
$wsBatting.Cells.Item(125, 4).Value = "4356"
$wsBatting.Cells.Item(126, 4).Value = "4357"
$wsBatting.Cells.Item(127, 4).Value = "4375"
$wsBatting.Cells.Item(128, 4).Value = "4376"

# Rows with no recorded INNING_NUMBER
$wsBatting.Cells.Item(3, 2).Value = $null
$wsBatting.Cells.Item(9, 2).Value = $null
$wsBatting.Cells.Item(12, 2).Value = $null
$wsBatting.Cells.Item(15, 2).Value = $null
$wsBatting.Cells.Item(21, 2).Value = $null
$wsBatting.Cells.Item(22, 2).Value = $null
$wsBatting.Cells.Item(26, 2).Value = $null
$wsBatting.Cells.Item(29, 2).Value = $null
$wsBatting.Cells.Item(30, 2).Value = $null
$wsBatting.Cells.Item(36, 2).Value = $null
$wsBatting.Cells.Item(38, 2).Value = $null
$wsBatting.Cells.Item(39, 2).Value = $null
$wsBatting.Cells.Item(42, 2).Value = $null
$wsBatting.Cells.Item(46, 2).Value = $null
$wsBatting.Cells.Item(55, 2).Value = $null
$wsBatting.Cells.Item(76, 2).Value = $null
$wsBatting.Cells.Item(80, 2).Value = $null
$wsBatting.Cells.Item(81, 2).Value = $null
$wsBatting.Cells.Item(84, 2).Value = $null
$wsBatting.Cells.Item(107, 2).Value = $null
$wsBatting.Cells.Item(126, 2).Value = $null

# --- "ODI Bowling": MATCH_CARD_LINK -> MATCH_CODE, URL -> bare match code ---
$wsBowling = $wb.Worksheets.Item("ODI Bowling")
$wsBowling.Cells.Item(1, 2).Value = "MATCH_CODE"
$wsBowling.Cells.Item(2, 2).Value = "3389"
$wsBowling.Cells.Item(3, 2).Value = "3527"
$wsBowling.Cells.Item(4, 2).Value = "3737"
$wsBowling.Cells.Item(5, 2).Value = "3739"

# --- New "Player Info" sheet, placed before "ODI Batting" ---
$wsPlayerInfo = $wb.Worksheets.Add($wsBatting)
$wsPlayerInfo.Name = "Player Info"
$wsPlayerInfo.Cells.Item(1, 1).Value = "ID"
$wsPlayerInfo.Cells.Item(1, 2).Value = "NAME"
$wsPlayerInfo.Cells.Item(1, 3).Value = "BATTING_HAND"
$wsPlayerInfo.Cells.Item(1, 4).Value = "BOWL_STYLE"
$wsPlayerInfo.Cells.Item(2, 1).Value = "3749"
$wsPlayerInfo.Cells.Item(2, 2).Value = "Hettige Don Rumesh Lahiru Thirimanne"
$wsPlayerInfo.Cells.Item(2, 3).Value = "Left Handed"
$wsPlayerInfo.Cells.Item(2, 4).Value = "Right Arm Medium Fast"
$piHeader = $wsPlayerInfo.Range("A1:D1")
$piHeader.Font.Bold = $true
$piHeader.Borders.LineStyle = 1
$piHeader.HorizontalAlignment = -4108
$piHeader.VerticalAlignment = -4160

# --- New "ODI Batting Extra" sheet, appended after "ODI Bowling" ---
$wsLast = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsExtra = $wb.Worksheets.Add($null, $wsLast)
$wsExtra.Name = "ODI Batting Extra"
$wsExtra.Cells.Item(1, 1).Value = "MATCH_CODE"
$wsExtra.Cells.Item(1, 2).Value = "BATTING_POSITION"
$wsExtra.Cells.Item(1, 3).Value = "NUM_4"
$wsExtra.Cells.Item(1, 4).Value = "NUM_6"
$wsExtra.Cells.Item(1, 5).Value = "PERCENT_RUNS_OF_TOTAL"
$wsExtra.Cells.Item(1, 6).Value = "MAN_OF_MATCH"
$wsExtra.Cells.Item(2, 1).Value = "4063"
$wsExtra.Cells.Item(2, 6).Value = "NO"
$wsExtra.Cells.Item(3, 1).Value = "4064"
$wsExtra.Cells.Item(3, 2).Value = 4
$wsExtra.Cells.Item(3, 3).Value = "1"
$wsExtra.Cells.Item(3, 4).Value = "1"
$wsExtra.Cells.Item(3, 5).Value = "8.70%"
$wsExtra.Cells.Item(3, 6).Value = "NO"
$wsExtra.Cells.Item(4, 1).Value = "4065"
$wsExtra.Cells.Item(4, 6).Value = "NO"
$wsExtra.Cells.Item(5, 1).Value = "4079"
$wsExtra.Cells.Item(5, 6).Value = "NO"
$wsExtra.Cells.Item(6, 1).Value = "4081"
$wsExtra.Cells.Item(6, 2).Value = 4
$wsExtra.Cells.Item(6, 3).Value = "0"
$wsExtra.Cells.Item(6, 4).Value = "0"
$wsExtra.Cells.Item(6, 5).Value = "6.42%"
$wsExtra.Cells.Item(6, 6).Value = "NO"
$wsExtra.Cells.Item(7, 1).Value = "4082"
$wsExtra.Cells.Item(7, 2).Value = 5
$wsExtra.Cells.Item(7, 3).Value = "0"
$wsExtra.Cells.Item(7, 4).Value = "0"
$wsExtra.Cells.Item(7, 5).Value = "13.46%"
$wsExtra.Cells.Item(7, 6).Value = "NO"
$wsExtra.Cells.Item(8, 1).Value = "4084"
$wsExtra.Cells.Item(8, 2).Value = 4
$wsExtra.Cells.Item(8, 3).Value = "4"
$wsExtra.Cells.Item(8, 4).Value = "0"
$wsExtra.Cells.Item(8, 5).Value = "35.84%"
$wsExtra.Cells.Item(8, 6).Value = "NO"
$wsExtra.Cells.Item(9, 1).Value = "4087"
$wsExtra.Cells.Item(9, 2).Value = 4
$wsExtra.Cells.Item(9, 3).Value = "3"
$wsExtra.Cells.Item(9, 4).Value = "0"
$wsExtra.Cells.Item(9, 5).Value = "18.45%"
$wsExtra.Cells.Item(9, 6).Value = "NO"
$wsExtra.Cells.Item(10, 1).Value = "4096"
$wsExtra.Cells.Item(10, 2).Value = 3
$wsExtra.Cells.Item(10, 3).Value = "0"
$wsExtra.Cells.Item(10, 4).Value = "0"
$wsExtra.Cells.Item(10, 6).Value = "NO"
$wsExtra.Cells.Item(11, 1).Value = "4098"
$wsExtra.Cells.Item(11, 6).Value = "NO"
$wsExtra.Cells.Item(12, 1).Value = "4302"
$wsExtra.Cells.Item(12, 6).Value = "NO"
$wsExtra.Cells.Item(13, 1).Value = "4305"
$wsExtra.Cells.Item(13, 2).Value = 1
$wsExtra.Cells.Item(13, 3).Value = "1"
$wsExtra.Cells.Item(13, 4).Value = "0"
$wsExtra.Cells.Item(13, 5).Value = "2.94%"
$wsExtra.Cells.Item(13, 6).Value = "NO"
$wsExtra.Cells.Item(14, 1).Value = "4309"
$wsExtra.Cells.Item(14, 2).Value = 3
$wsExtra.Cells.Item(14, 3).Value = "1"
$wsExtra.Cells.Item(14, 4).Value = "0"
$wsExtra.Cells.Item(14, 5).Value = "12.44%"
$wsExtra.Cells.Item(14, 6).Value = "NO"
$wsExtra.Cells.Item(15, 1).Value = "4322"
$wsExtra.Cells.Item(15, 2).Value = 3
$wsExtra.Cells.Item(15, 3).Value = "1"
$wsExtra.Cells.Item(15, 4).Value = "0"
$wsExtra.Cells.Item(15, 5).Value = "6.48%"
$wsExtra.Cells.Item(15, 6).Value = "NO"
$wsExtra.Cells.Item(16, 1).Value = "4344"
$wsExtra.Cells.Item(16, 2).Value = 6
$wsExtra.Cells.Item(16, 3).Value = "4"
$wsExtra.Cells.Item(16, 4).Value = "0"
$wsExtra.Cells.Item(16, 5).Value = "13.31%"
$wsExtra.Cells.Item(16, 6).Value = "NO"
$wsExtra.Cells.Item(17, 1).Value = "4350"
$wsExtra.Cells.Item(17, 2).Value = 6
$wsExtra.Cells.Item(17, 3).Value = "4"
$wsExtra.Cells.Item(17, 4).Value = "0"
$wsExtra.Cells.Item(17, 5).Value = "20.08%"
$wsExtra.Cells.Item(17, 6).Value = "NO"
$wsExtra.Cells.Item(18, 1).Value = "4356"
$wsExtra.Cells.Item(18, 6).Value = "NO"
$wsExtra.Cells.Item(19, 1).Value = "4357"
$wsExtra.Cells.Item(19, 6).Value = "NO"
$wsExtra.Cells.Item(20, 1).Value = "4375"
$wsExtra.Cells.Item(20, 2).Value = 5
$wsExtra.Cells.Item(20, 3).Value = "0"
$wsExtra.Cells.Item(20, 4).Value = "0"
$wsExtra.Cells.Item(20, 6).Value = "NO"
$wsExtra.Cells.Item(21, 1).Value = "4376"
$wsExtra.Cells.Item(21, 2).Value = 3
$wsExtra.Cells.Item(21, 3).Value = "4"
$wsExtra.Cells.Item(21, 4).Value = "0"
$wsExtra.Cells.Item(21, 5).Value = "12.12%"
$wsExtra.Cells.Item(21, 6).Value = "NO"
$exHeader = $wsExtra.Range("A1:F1")
$exHeader.Font.Bold = $true
$exHeader.Borders.LineStyle = 1
$exHeader.HorizontalAlignment = -4108
$exHeader.VerticalAlignment = -4160

